{"js": "// Update the Solicita\u00e7\u00e3o de Compartilhamento document with the new\n// requester's (SNI NET TELECOM) data: city/state/date, company name,\n// CNPJ, address, signatory name/title, e-mail and phone.\nconst replacements = [\n  [\"Parnaiba\", \"Mat\u00f5es\"],\n  [\"Piau\u00ed\", \"Maranh\u00e3o\"],\n  [\"30 de janeiro de 2026\", \"4 de fevereiro de 2026\"],\n  [\"EQUATORIAL PIAU\u00cd\", \"EQUATORIAL MARANH\u00c3O\"],\n  [\"INFORLAN\", \"SNI NET TELECOM\"],\n  [\"SISTEMA INFOR-LAN TELECOMUNICA\u00c7\u00d5ES LTDA\", \"C. O. DOS SANTOS SOUSA\"],\n  [\"52.629.625/0001-10\", \"06.323.714/0001-17\"],\n  [\n    \"Q QUADRA 1, 08, PLANALTO DE MONTESERRA THE, Parna\u00edba/PI, CEP: 64.207-470\",\n    \"R Timon, 355, Centro, Parnarama/MA, CEP: 65.640-000\",\n  ],\n  [\"MARCIA MARIA MENDES COSTA\", \"Carlos Odeon dos Santos Sousa\"],\n  [\"S\u00f3cia-Administradora\", \"S\u00f3cio/Propri\u00e9tario\"],\n  [\"sistemainforlan@gmail.com\", \"CARLOS16SNI@HOTMAIL.COM\"],\n  [\"(86) 9428-1942\", \"(99) 8430-1702\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the company/address/date/signature details for the new\n# requester (SNI NET TELECOM) per the commit.\n$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n\nReplace-Text \"Parnaiba\" \"Mat\u00f5es\"\nReplace-Text \"Piau\u00ed\" \"Maranh\u00e3o\"\nReplace-Text \"30 de janeiro de 2026\" \"4 de fevereiro de 2026\"\nReplace-Text \"EQUATORIAL PIAU\u00cd\" \"EQUATORIAL MARANH\u00c3O\"\nReplace-Text \"INFORLAN\" \"SNI NET TELECOM\"\nReplace-Text \"SISTEMA INFOR-LAN TELECOMUNICA\u00c7\u00d5ES LTDA\" \"C. O. DOS SANTOS SOUSA\"\nReplace-Text \"52.629.625/0001-10\" \"06.323.714/0001-17\"\nReplace-Text \"Q QUADRA 1, 08, PLANALTO DE MONTESERRA THE, Parna\u00edba/PI, CEP: 64.207-470\" \"R Timon, 355, Centro, Parnarama/MA, CEP: 65.640-000\"\nReplace-Text \"MARCIA MARIA MENDES COSTA\" \"Carlos Odeon dos Santos Sousa\"\nReplace-Text \"S\u00f3cia-Administradora\" \"S\u00f3cio/Propri\u00e9tario\"\nReplace-Text \"sistemainforlan@gmail.com\" \"CARLOS16SNI@HOTMAIL.COM\"\nReplace-Text \"(86) 9428-1942\" \"(99) 8430-1702\"\n"}
